$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: B1:Z1 (x values)
$row1 = @(9.21771, 4.25799, 4.42522, 0.397115, 2.8813, 4.87695, 7.91005, 7.70865, 6.85066, 4.51345, 8.03122, 5.39499, 7.01773, 5.879, 5.87789, 6.29149, 4.7437, 1.62803, 8.55286, 3.77229, 5.93972, 3.8352, 6.33699, 4.79261, 5.29877)
for ($i = 0; $i -lt $row1.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $row1[$i]
}

# Row 2: B2:Z2 (y values)
$row2 = @(-0.47033, -0.23693, -1.4953, -0.22986, -0.42457, -1.559, 0.15847, 0.52688, -0.0065266, 1.0268, 0.16051, -2.3574, -2.6056, 0.033404, -0.58253, 1.7387, -0.98227, 0.94823, -2.6024, -0.14505, 0.36789, 0.32945, 0.30675, -0.0080439, -0.41344)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 2).Value = $row2[$i]
}

# F11: new literal value
$ws.Range("F11").Value = -21.365194937111

# Row 18: A18:AX18
$row18 = @(4.70724, 5.26906, 1.25723, 5.9213, 4.57149, 5.32742, 3.34411, 5.59719, 7.11093, 5.02043, 7.34914, 3.90632, 2.90111, 6.32136, 3.74945, 7.97193, 3.34184, -0.118245, 3.22259, 3.92044, 7.03845, 3.74209, 4.03482, 5.67917, 4.75739, 9.21771, 4.25799, 4.42522, 0.397115, 2.8813, 4.87695, 7.91005, 7.70865, 6.85066, 4.51345, 8.03122, 5.39499, 7.01773, 5.879, 5.87789, 6.29149, 4.7437, 1.62803, 8.55286, 3.77229, 5.93972, 3.8352, 6.33699, 4.79261, 5.29877)
for ($i = 0; $i -lt $row18.Length; $i++) {
    $ws.Cells.Item(18, $i + 1).Value = $row18[$i]
}

# Row 19: A19:Y19
$row19 = @(-0.47033, -0.23693, -1.4953, -0.22986, -0.42457, -1.559, 0.15847, 0.52688, -0.0065266, 1.0268, 0.16051, -2.3574, -2.6056, 0.033404, -0.58253, 1.7387, -0.98227, 0.94823, -2.6024, -0.14505, 0.36789, 0.32945, 0.30675, -0.0080439, -0.41344)
for ($i = 0; $i -lt $row19.Length; $i++) {
    $ws.Cells.Item(19, $i + 1).Value = $row19[$i]
}

# Selection
$ws.Range("K14").Select()

